# Edit script: rewrite "quantum realm" themed document into "chemistry" themed
# document, update author/byline, and set the font from the misspelled
# "TimesNewToman" to the correct "Times New Roman".

$d = $word.ActiveDocument

function Replace-Text {
    param(
        [string]$Find,
        [string]$Replacement
    )
    $d.Content.Find.Execute($Find, $true, $false, $false, $false, $false, $true, 1, $false, $Replacement, 2) | Out-Null
}

# ---------------------------------------------------------------------------
# 1. Fix the font across the whole document (TimesNewToman -> Times New Roman)
# ---------------------------------------------------------------------------
$fullRange = $d.Range(0, $d.Content.End)
$fullRange.Font.Name = "Times New Roman"

# ---------------------------------------------------------------------------
# 2. Restructure the long "quantum" paragraph's tail:
#    " These discoveries ... quantum mechanics.The study of the quantum
#    realm is not merely ... navigation." gets replaced by two manual line
#    breaks followed by two new chemistry-themed sentences.
#    (Find text is limited to ~255 chars, so we locate short anchors and
#    then rewrite the whole Range between them.)
# ---------------------------------------------------------------------------
$anchorStart = $d.Content
$anchorStart.Find.Execute(" These discoveries have opened up", $true, $false, $false, $false, $false, $true, 1, $false) | Out-Null
$spanStart = $anchorStart.Start

$anchorEnd = $d.Content
$anchorEnd.Find.Execute("materials science, and navigation.", $true, $false, $false, $false, $false, $true, 1, $false) | Out-Null
$spanEnd = $anchorEnd.End

$manualBreak = [string][char]11
$newTail = $manualBreak + $manualBreak + "Chemistry opens a vast canvas of opportunities for exploration and discovery, empowering us to unravel the enigmas of the microscopic world and pushing the boundaries of scientific understanding. Whether it's deciphering the intricate workings of cells, developing innovative materials with remarkable properties, or devising groundbreaking medical treatments, chemistry serves as a catalyst for progress and transformation in countless fields of human endeavor."

$tailRange = $d.Range($spanStart, $spanEnd)
$tailRange.Text = $newTail

# ---------------------------------------------------------------------------
# 3. Title / byline / contact info
# ---------------------------------------------------------------------------
Replace-Text "Quantum Realm Unveiled: The Nexus of Science and Imagination" "The Marvelous Machine: Exploring Chemistry and Its Role in Our Lives"
Replace-Text "Dr" "Beatrice A"
Replace-Text " Alex Hayes" " Franklin"
Replace-Text "hayes" "franklin"
Replace-Text "alex@academic" "bea@schoolmail"
Replace-Text "edu" "org"

# ---------------------------------------------------------------------------
# 4. Body paragraph sentences
# ---------------------------------------------------------------------------
Replace-Text "From the dawn of time, humanity has pondered the enigmatic realm of the quantum world, a domain where particles behave in ways that defy classical intuition" "Our world is an intricate tapestry woven together by countless chemical reactions, each contributing to the vibrant spectacle of life"
Replace-Text " This realm lies at the heart of matter, where subatomic particles dance in a cosmic symphony of probability and indeterminacy" " Chemistry, the study of matter and its properties, offers a magnifying glass into these intricate processes, revealing the fundamental building blocks of our universe and their interactions"
Replace-Text " In this ethereal realm, physicists seek to unravel the fundamental laws that govern the universe, pushing the boundaries of human knowledge" " From the air we breathe, to the food we eat, and the medicines that heal us, chemistry is an omnipresent force, shaping our world in myriad ways"
Replace-Text "Delving into the quantum realm has led to profound insights into the nature of reality" "As we delve into the realm of chemistry, we uncover a fascinating dance of atoms and molecules, a symphony of interactions governed by intricate laws"
Replace-Text " Experiments have revealed that particles can exist in multiple states simultaneously, defying our everyday notions of locality" " From towering mountains sculpted by weathering to the burning of a simple candle, the principles of chemistry intricately orchestrate the countless phenomena that unfold around us"
Replace-Text " The phenomenon of entanglement, where particles separated by vast distances remain mysteriously interconnected, challenges our understanding of cause and effect" " It is through chemistry that we can comprehend the mysteries of the natural world, unlocking its secrets and harnessing its power for human benefit"

# ---------------------------------------------------------------------------
# 5. Summary paragraph sentences
# ---------------------------------------------------------------------------
Replace-Text "Our exploration of the quantum realm has revealed a universe governed by laws far removed from our everyday experience" "The study of chemistry reveals the fundamental building blocks of our universe and their interactions, offering insights into the myriad chemical reactions that shape our world"
Replace-Text " Quantum mechanics has challenged our understanding of reality, revealing the strange and wondrous behaviors of subatomic particles" " Chemistry empowers us to comprehend the mysteries of the natural world, unlocking its secrets and harnessing its power for human benefit"
Replace-Text " This newfound knowledge promises to fuel technological revolutions, enabling breakthroughs in computing, cryptography, and sensing" " It opens up vast avenues of exploration and discovery, enabling us to decipher the inner workings of cells, develop innovative materials, and devise life-saving treatments"
Replace-Text " As we continue to unravel the mysteries of the quantum realm, we stand at the threshold of an era where science and imagination converge, transforming our understanding of the universe and reshaping the world we live in" " Chemistry stands as a testament to the interconnectedness of all matter, underscoring the profound impact it has on our lives and the world around us"

# ---------------------------------------------------------------------------
# 6. Add a trailing empty paragraph at the end of the document
# ---------------------------------------------------------------------------
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()

# ---------------------------------------------------------------------------
# 7. Make sure the whole document (including newly-inserted text) uses the
#    corrected font name.
# ---------------------------------------------------------------------------
$fullRange = $d.Range(0, $d.Content.End)
$fullRange.Font.Name = "Times New Roman"
